$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Daniel M.%Bean%NULL%0,                  Zeljko%Kraljevic%NULL%1,                  Thomas%Searle%NULL%1,                  Rebecca%Bendayan%NULL%2,                  O''Gallagher%Kevin%NULL%1,                  Andrew%Pickles%NULL%1,                  Amos%Folarin%NULL%1,                  Lukasz%Roguski%NULL%1,                  Kawsar%Noor%NULL%1,                  Anthony%Shek%NULL%1,                  Rosita%Zakeri%NULL%1,                  Ajay M.%Shah%ajay.shah@kcl.ac.uk%2,                  James T.H.%Teo%jamesteo@nhs.net%2,                  Richard J.B.%Dobson%richard.j.dobson@kcl.ac.uk%2]'
$ws.Range("I2").Value = '_PMC_CROSSREF'

$ws.Range("E3").Value = '[Harmony R.%Reynolds%NULL%0,                  Samrachana%Adhikari%NULL%1,                  Claudia%Pulgarin%NULL%1,                  Andrea B.%Troxel%NULL%1,                  Eduardo%Iturrate%NULL%1,                  Stephen B.%Johnson%NULL%1,                  Anaïs%Hausvater%NULL%1,                  Jonathan D.%Newman%NULL%1,                  Jeffrey S.%Berger%NULL%1,                  Sripal%Bangalore%NULL%1,                  Stuart D.%Katz%NULL%1,                  Glenn I.%Fishman%NULL%1,                  Dennis%Kunichoff%NULL%1,                  Yu%Chen%NULL%1,                  Gbenga%Ogedegbe%NULL%1,                  Judith S.%Hochman%NULL%1]'
$ws.Range("I3").Value = '_PMC_CROSSREF'

$ws.Range("E4").Value = '[Giuseppe%Mancia%NULL%0,                  Federico%Rea%NULL%1,                  Monica%Ludergnani%NULL%2,                  Monica%Ludergnani%NULL%0,                  Giovanni%Apolone%NULL%1,                  Giovanni%Corrao%NULL%1]'
$ws.Range("I4").Value = '_PMC_CROSSREF'

$ws.Range("E5").Value = '[Mandeep R.%Mehra%NULL%1,                  Sapan S.%Desai%NULL%2,                  Sapan S.%Desai%NULL%0,                  SreyRam%Kuy%NULL%2,                  SreyRam%Kuy%NULL%0,                  Timothy D.%Henry%NULL%1,                  Amit N.%Patel%NULL%1]'
$ws.Range("I5").Value = '_PMC_CROSSREF'

$ws.Range("E6").Value = '[ Z.%Feng%null%1,                   J.% Li%null%1,                   S.% Yao%null%1,                   Q.% Yu%null%1,                   W.% Zhou%null%1,                   X.% Mao%null%1,                   H.% Li%null%1,                   W.% Kang%null%1,                   X.% Ouyang%null%1,                   J.% Mei%null%1,                   Q.% Zeng%null%1,                   J.% Liu%null%1,                   X.% Ma%null%1,                   P.% Rong%null%1,                   W. % Wang%null%1]'
$ws.Range("I6").Value = '_MedBiorxiv_CROSSREF'

$ws.Range("E7").Value = '[Guang%Yang%xref no email%0,            Zihu%Tan%xref no email%0,            Ling%Zhou%xref no email%0,            Min%Yang%xref no email%0,            Lang%Peng%xref no email%0,            Jinjin%Liu%xref no email%0,            Jingling%Cai%xref no email%0,            Ru%Yang%xref no email%0,            Junyan%Han%xref no email%0,            Yafei%Huang%xref no email%0,            Shaobin%He%xref no email%0]'
$ws.Range("I7").Value = '_CROSSREF'

$ws.Range("E8").Value = '[Peng%Zhang%NULL%0,                  Lihua%Zhu%NULL%1,                  Jingjing%Cai%NULL%1,                  Fang%Lei%NULL%1,                  Juan-Juan%Qin%NULL%1,                  Jing%Xie%NULL%0,                  Ye-Mao%Liu%NULL%1,                  Yan-Ci%Zhao%NULL%1,                  Xuewei%Huang%NULL%1,                  Lijin%Lin%NULL%1,                  Meng%Xia%NULL%1,                  Ming-Ming%Chen%NULL%1,                  Xu%Cheng%NULL%1,                  Xiao%Zhang%NULL%1,                  Deliang%Guo%NULL%1,                  Yuanyuan%Peng%NULL%1,                  Yan-Xiao%Ji%NULL%1,                  Jing%Chen%NULL%0,                  Zhi-Gang%She%NULL%1,                  Yibin%Wang%NULL%1,                  Qingbo%Xu%NULL%1,                  Renfu%Tan%NULL%1,                  Haitao%Wang%NULL%1,                  Jun%Lin%NULL%1,                  Pengcheng%Luo%NULL%1,                  Shouzhi%Fu%NULL%2,                  Hongbin%Cai%NULL%1,                  Ping%Ye%NULL%1,                  Bing%Xiao%NULL%1,                  Weiming%Mao%NULL%1,                  Liming%Liu%NULL%1,                  Youqin%Yan%NULL%1,                  Mingyu%Liu%NULL%1,                  Manhua%Chen%NULL%1,                  Xiao-Jing%Zhang%NULL%1,                  Xinghuan%Wang%NULL%2,                  Rhian M.%Touyz%NULL%1,                  Jiahong%Xia%NULL%1,                  Bing-Hong%Zhang%NULL%1,                  Xiaodong%Huang%NULL%1,                  Yufeng%Yuan%NULL%1,                  Rohit%Loomba%NULL%1,                  Peter P.%Liu%NULL%1,                  Hongliang%Li%NULL%1]'
$ws.Range("I8").Value = '_PMC_CROSSREF'

$ws.Range("E9").Value = '[Juan%Meng%NULL%0,                  Guohui%Xiao%NULL%1,                  Juanjuan%Zhang%NULL%1,                  Xing%He%NULL%1,                  Min%Ou%NULL%1,                  Jing%Bi%NULL%1,                  Rongqing%Yang%NULL%1,                  Wencheng%Di%NULL%1,                  Zhaoqin%Wang%NULL%2,                  Zigang%Li%NULL%1,                  Hong%Gao%NULL%0,                  Lei%Liu%NULL%0,                  Guoliang%Zhang%NULL%1]'
$ws.Range("I9").Value = '_PMC_CROSSREF'

$ws.Range("E10").Value = '[Lei-Ke%Zhang%NULL%1,                  Yuan%Sun%NULL%2,                  Yuan%Sun%NULL%0,                  Haolong%Zeng%NULL%2,                  Haolong%Zeng%NULL%0,                  Qingxing%Wang%NULL%1,                  Xiaming%Jiang%NULL%1,                  Wei-Juan%Shang%NULL%1,                  Yan%Wu%NULL%1,                  Shufen%Li%NULL%1,                  Yu-Lan%Zhang%NULL%1,                  Zhao-Nian%Hao%NULL%1,                  Hongbo%Chen%NULL%1,                  Runming%Jin%NULL%2,                  Runming%Jin%NULL%0,                  Wei%Liu%NULL%0,                  Hao%Li%lihao_1986@126.com%1,                  Ke%Peng%pengke@wh.iov.cn%1,                  Gengfu%Xiao%xiaogf@wh.iov.cn%2,                  Gengfu%Xiao%xiaogf@wh.iov.cn%0]'
$ws.Range("I10").Value = '_PMC_Springer_CROSSREF'

$ws.Range("E11").Value = '[Renhong%Yan%xref no email%1,            Yuanyuan%Zhang%xref no email%1,            Yaning%Li%xref no email%1,            Lu%Xia%xref no email%1,            Qiang%Zhou%xref no email%1]'
$ws.Range("I11").Value = '_CROSSREF'

$ws.Range("E12").Value = '[Juyi%Li%xref no email%0,            Xiufang%Wang%xref no email%0,            Jian%Chen%xref no email%0,            Hongmei%Zhang%xref no email%0,            Aiping%Deng%xref no email%0]'
$ws.Range("I12").Value = '_CROSSREF'

$ws.Range("E13").Value = '[Tao%Guo%NULL%5,                  Yongzhen%Fan%NULL%1,                  Ming%Chen%NULL%1,                  Xiaoyan%Wu%NULL%1,                  Lin%Zhang%NULL%1,                  Tao%He%NULL%1,                  Hairong%Wang%NULL%1,                  Jing%Wan%NULL%1,                  Xinghuan%Wang%NULL%0,                  Zhibing%Lu%NULL%1]'
$ws.Range("I13").Value = '_PMC_CROSSREF'
